# Atualizado por script em 11-11-2023 08:45
#
# The source site re-scraped the Ligue 1 2023-2024 schedule: three rows in
# the already-present block (86-88, 94-96) got cyclically re-ordered
# (same matches, different row position) and one brand-new match
# (Montpellier x Nice, played 10/11/2023) was appended as row 99.
#
# NOTE: this interpreter's functions only bind parameters positionally
# (named "-Param value" binding does not work), so Set-MatchRow takes its
# arguments in a fixed order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-MatchRow($Row, $Home, $HomeGoals, $Away, $AwayGoals, `
        $HomeOpenOdds, $HomeOpenDt, $HomeCloseOdds, $HomeCloseDt, `
        $DrawOpenOdds, $DrawOpenDt, $DrawCloseOdds, $DrawCloseDt, `
        $AwayOpenOdds, $AwayOpenDt, $AwayCloseOdds, $AwayCloseDt, $Url) {

    $ws.Cells.Item($Row, 6).Value = $Home
    $ws.Cells.Item($Row, 7).Value = $HomeGoals
    $ws.Cells.Item($Row, 8).Value = $Away
    $ws.Cells.Item($Row, 9).Value = $AwayGoals
    $ws.Cells.Item($Row, 10).Value = $HomeOpenOdds
    $ws.Cells.Item($Row, 11).Value = $HomeOpenDt
    $ws.Cells.Item($Row, 12).Value = $HomeCloseOdds
    $ws.Cells.Item($Row, 13).Value = $HomeCloseDt
    $ws.Cells.Item($Row, 14).Value = $DrawOpenOdds
    $ws.Cells.Item($Row, 15).Value = $DrawOpenDt
    $ws.Cells.Item($Row, 16).Value = $DrawCloseOdds
    $ws.Cells.Item($Row, 17).Value = $DrawCloseDt
    $ws.Cells.Item($Row, 18).Value = $AwayOpenOdds
    $ws.Cells.Item($Row, 19).Value = $AwayOpenDt
    $ws.Cells.Item($Row, 20).Value = $AwayCloseOdds
    $ws.Cells.Item($Row, 21).Value = $AwayCloseDt
    $ws.Cells.Item($Row, 22).Value = $Url
}

# --- Row 86 becomes "Metz vs Le Havre" (previously row 88) ---
Set-MatchRow 86 "Metz" 0 "Le Havre" 0 `
    2.43 "11/10/2023 14:10" 2.95 "29/10/2023 14:58" `
    3.32 "11/10/2023 14:10" 3.09 "29/10/2023 14:31" `
    2.92 "11/10/2023 14:10" 2.74 "29/10/2023 14:58" `
    "https://www.betexplorer.com/football/france/ligue-1/metz-le-havre/t0JqDX2r/"

# --- Row 87 becomes "Montpellier vs Toulouse" (previously row 86) ---
Set-MatchRow 87 "Montpellier" 3 "Toulouse" 0 `
    1.67 "10/10/2023 14:02" 1.88 "29/10/2023 14:59" `
    4.05 "10/10/2023 14:02" 3.82 "29/10/2023 14:59" `
    4.64 "10/10/2023 14:02" 4.24 "29/10/2023 14:59" `
    "https://www.betexplorer.com/football/france/ligue-1/montpellier-toulouse/KjImCiIl/"

# --- Row 88 becomes "Lille vs Monaco" (previously row 87) ---
Set-MatchRow 88 "Lille" 2 "Monaco" 0 `
    2.01 "10/10/2023 14:02" 2.5 "29/10/2023 14:52" `
    3.81 "10/10/2023 14:02" 3.7 "29/10/2023 14:52" `
    3.62 "10/10/2023 14:02" 2.79 "29/10/2023 14:52" `
    "https://www.betexplorer.com/football/france/ligue-1/lille-monaco/l2v5KFA8/"

# --- Row 94 becomes "Strasbourg vs Clermont" (previously row 95) ---
Set-MatchRow 94 "Strasbourg" 0 "Clermont" 0 `
    2 "22/10/2023 12:02" 2.33 "05/11/2023 14:59" `
    3.45 "22/10/2023 12:02" 3.23 "05/11/2023 14:59" `
    4.06 "22/10/2023 12:02" 3.46 "05/11/2023 14:59" `
    "https://www.betexplorer.com/football/france/ligue-1/strasbourg-clermont/YFXh9k27/"

# --- Row 95 becomes "Toulouse vs Le Havre" (previously row 96) ---
Set-MatchRow 95 "Toulouse" 1 "Le Havre" 2 `
    2.17 "23/10/2023 15:49" 2.04 "05/11/2023 14:36" `
    3.41 "23/10/2023 15:49" 3.44 "05/11/2023 14:39" `
    3.56 "23/10/2023 15:49" 4.05 "05/11/2023 14:36" `
    "https://www.betexplorer.com/football/france/ligue-1/toulouse-le-havre/nJyc89HD/"

# --- Row 96 becomes "Nantes vs Reims" (previously row 94) ---
Set-MatchRow 96 "Nantes" 0 "Reims" 1 `
    2.87 "22/10/2023 12:02" 2.93 "05/11/2023 14:54" `
    3.44 "22/10/2023 12:02" 3.3 "05/11/2023 14:58" `
    2.41 "22/10/2023 12:02" 2.61 "05/11/2023 14:53" `
    "https://www.betexplorer.com/football/france/ligue-1/nantes-reims/fPpM8wTg/"

# --- New row 99: Montpellier vs Nice ---
# Columns A (Indice) and E (data_partida) carry dedicated styles in this
# sheet (bold/centered/bordered index column, datetime-formatted date
# column); copy those formats from the row right above rather than
# hand-building them, so no stray style entries get created.
$ws.Range("A98").Copy() | Out-Null
$ws.Range("A99").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E98").Copy() | Out-Null
$ws.Range("E99").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(99, 1).Value = 98
$ws.Cells.Item(99, 2).Value = "france"
$ws.Cells.Item(99, 3).Value = "ligue-1"
$ws.Cells.Item(99, 4).Value = "2023-2024"
$ws.Cells.Item(99, 5).Value = 45240.875

Set-MatchRow 99 "Montpellier" 0 "Nice" 0 `
    2.74 "29/10/2023 11:02" 3.69 "10/11/2023 20:59" `
    3.32 "29/10/2023 11:02" 3.36 "10/11/2023 20:51" `
    2.68 "29/10/2023 11:02" 2.18 "10/11/2023 20:53" `
    "https://www.betexplorer.com/football/france/ligue-1/montpellier-nice/zq6eNxrm/"
